$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: write the new / changed numeric values in column H (and F10:F12) ---
$ws.Range("H1").Value = 6
$ws.Range("H2").Value = 12

$ws.Range("H3").Value = 6
$ws.Range("H4").Value = 5
$ws.Range("H5").Value = 5

$ws.Range("H6").Value = 2
$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 7

$ws.Range("H9").Value = 1

$ws.Range("F10").Value = 1200
$ws.Range("H10").Value = 1

$ws.Range("F11").Value = 1200
$ws.Range("H11").Value = 5

$ws.Range("F12").Value = 1200
$ws.Range("H12").Value = 4

$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 3
$ws.Range("H15").Value = 3
$ws.Range("H16").Value = 3

# --- Step 2: mark the same rows with a red font (applied last so the new values keep the style) ---
$ws.Range("B6:H6").Font.Color = 255
$ws.Range("B7:H7").Font.Color = 255
$ws.Range("B9:H9").Font.Color = 255

$ws.Range("B10:F10").Font.Color = 255
$ws.Range("B11:F11").Font.Color = 255
$ws.Range("B12:F12").Font.Color = 255

$ws.Range("B13:F13").Font.Color = 255
$ws.Range("H13").Font.Color = 255
$ws.Range("B14:F14").Font.Color = 255
$ws.Range("H14").Font.Color = 255
$ws.Range("B15:F15").Font.Color = 255
$ws.Range("H15").Font.Color = 255
$ws.Range("B16:F16").Font.Color = 255
$ws.Range("H16").Font.Color = 255

# --- Step 3: move the active selection ---
$null = $ws.Range("D9").Select()

# --- Step 4: page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = $xlPortrait
